$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing row 3 ("Semester 2", ...) down to row 4
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new "Semester 1" entry.
# B3 and E3 must be stored as text (matching "2" and "4.00"), not as numbers,
# so force the Text number format before assigning those values.
$ws.Range("A3").Value = "Semester 1"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2"

$ws.Range("C3").Value = 2

$ws.Range("D3").Value = "A"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.00"
